$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "30.602.43"
$ws.Cells.Item(2, 5).Value = "  +0.62%  "

$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "1.923.49"
$ws.Cells.Item(3, 5).Value = "  -0.19%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  +0.15%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "247.36"
$ws.Cells.Item(5, 5).Value = "  +2.79%  "

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "1.001"
$ws.Cells.Item(6, 5).Value = "  +0.09%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4746"
$ws.Cells.Item(7, 5).Value = "  -0.13%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.2900"
$ws.Cells.Item(8, 5).Value = "  +1.14%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.06819"
$ws.Cells.Item(9, 5).Value = "  +3.90%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "105.27"
$ws.Cells.Item(10, 5).Value = "  -1.62%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "18.40"
$ws.Cells.Item(11, 5).Value = "  -3.36%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "1.927.15"
$ws.Cells.Item(12, 5).Value = "  +0.28%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.07698"
$ws.Cells.Item(13, 5).Value = "  +1.08%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.356"
$ws.Cells.Item(14, 5).Value = "  +4.45%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.6698"
$ws.Cells.Item(15, 5).Value = "  +2.12%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "290.58"
$ws.Cells.Item(16, 5).Value = "  -4.94%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "30.611.61"
$ws.Cells.Item(17, 5).Value = "  +0.63%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000007615"
$ws.Cells.Item(18, 5).Value = "  +1.75%  "

$ws.Cells.Item(19, 5).Value = "  +0.00%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "12.95"
$ws.Cells.Item(20, 5).Value = "  +0.04%  "

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "5.553"
$ws.Cells.Item(21, 5).Value = "  +4.87%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "2.182.31"
$ws.Cells.Item(22, 5).Value = "  +0.21%  "

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "1.001"
$ws.Cells.Item(23, 5).Value = "  +0.13%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "6.460"
$ws.Cells.Item(24, 5).Value = "  +3.28%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "9.520"
$ws.Cells.Item(25, 5).Value = "  +3.36%  "

$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "167.68"
$ws.Cells.Item(26, 5).Value = "  +0.22%  "

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "20.90"
$ws.Cells.Item(27, 5).Value = "  +3.75%  "

$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.122"
$ws.Cells.Item(28, 5).Value = "  +5.02%  "

$ws.Cells.Item(29, 5).Value = "  -3.34%  "

$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.404"
$ws.Cells.Item(30, 5).Value = "  +3.56%  "

$ws.Cells.Item(31, 5).Value = "  +2.42%  "

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "4.047"
$ws.Cells.Item(32, 5).Value = "  +3.41%  "

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "0.05019"
$ws.Cells.Item(33, 5).Value = "  +0.54%  "

$ws.Cells.Item(34, 5).Value = "  -0.95%  "

$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "1.145"
$ws.Cells.Item(35, 5).Value = "  -0.10%  "

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "0.02069"
$ws.Cells.Item(36, 5).Value = "  +6.74%  "

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.9995"
$ws.Cells.Item(37, 5).Value = "  +0.00%  "

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "2.722"
$ws.Cells.Item(38, 5).Value = "  -1.02%  "

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.687"
$ws.Cells.Item(39, 5).Value = "  -0.34%  "

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "111.67"
$ws.Cells.Item(40, 5).Value = "  +4.53%  "

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "2.037"
$ws.Cells.Item(41, 5).Value = "  -0.74%  "

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.8743"
$ws.Cells.Item(42, 5).Value = "  -0.38%  "

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.4407"
$ws.Cells.Item(43, 5).Value = "  +6.48%  "

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "5.902"
$ws.Cells.Item(44, 5).Value = "  +1.86%  "

$ws.Cells.Item(45, 5).Value = "  +0.09%  "

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "67.76"
$ws.Cells.Item(46, 5).Value = "  -3.19%  "

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "7.287"
$ws.Cells.Item(47, 5).Value = "  +0.45%  "

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "9.341"
$ws.Cells.Item(48, 5).Value = "  +0.75%  "

$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "48.20"
$ws.Cells.Item(49, 5).Value = "  +14.05%  "

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.1243"
$ws.Cells.Item(50, 5).Value = "  +3.52%  "

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "35.01"
$ws.Cells.Item(51, 5).Value = "  +0.62%  "
